$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C35").Value = "world reorganization & data workflow update"
$ws.Range("G35").Value = 3

$ws.Range("G39").Formula = "=SUM(G4:G35)"

$ws.Range("C36").Select()
$excel.ActiveWindow.ScrollRow = 25

